# Case_Data.xlsx - "Addressed all JTC scenarios in code and template."
#
# 1. Remove the two stray empty inlineStr cells at G1168 / G1169.
# 2. Append 17 new docket rows (1171-1187) for case 21CRB01268 / Hemmeter /
#    "Possession Drug Paraphernalia" (JTC = Juvenile Traffic Court? -> just
#    more charge rows being logged), which pushes the sheet's used range
#    from A2:K1170 to A2:K1187.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Drop the orphan empty-string cells in column G -------------------
$ws.Cells.Item(1168, 7).ClearContents()
$ws.Cells.Item(1169, 7).ClearContents()

# --- 2. Append the new rows ------------------------------------------------
# Columns: A Case#, B Defendant, C Charge, D Statute, E Class, F Plea,
#          G Finding, H Fine, I Costs, J Days, K Suspended
$newRows = @(
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","3","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","2"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","2"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","2"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","3"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","3"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","3"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None"),
    @("21CRB01268","Hemmeter","Possession Drug Paraphernalia","2925.14(C)","M4","No Contest","Guilty","$ 50","$ 0","5","None")
)

$startRow = 1171
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $text = $vals[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        # Values that Excel would otherwise auto-coerce (currency-looking
        # "$ 50" or plain digit strings like "3"/"5"/"2") need the cell
        # pre-formatted as Text so they are stored as literal strings,
        # exactly like the rest of this column in the existing data.
        if ($text -match '^\$ ?\d+(\.\d+)?$' -or $text -match '^\d+(\.\d+)?$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $text
    }
}
